# Update grading-comment cells (column F) on the "CustomerMappingDriver Class"
# and "Generic" sections of the rubric sheet, reflecting a closer review of
# driver-class code (lines 65-80).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F29").Value = "(-8) points for not initializing Customer, not checking condition correctly, not declaring and initializing product and adding them to inventory"
$ws.Range("F37").Value = "(-5) for compilation errors in CustomerMapping class"
$ws.Range("F30").Value = "(-4) for no output due to Compilation errors"

# Reflect where the grader's cursor ended up after making the edits.
$ws.Range("F30").Select()
